$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows 5-7 replicate existing rows 3 (Dubai), 2 (Abu Dhabi), 4 (Sharjah).
# Copy/PasteSpecial (rather than re-typing values) keeps every cell's original
# data type (the numeric-looking columns G:K stay stored as text) without
# introducing new cell styles.
$xlPasteAll = -4104

$ws.Range("A3:K3").Copy()
$ws.Range("A5:K5").PasteSpecial($xlPasteAll)

$ws.Range("A2:K2").Copy()
$ws.Range("A6:K6").PasteSpecial($xlPasteAll)

$ws.Range("A4:K4").Copy()
$ws.Range("A7:K7").PasteSpecial($xlPasteAll)

$excel.CutCopyMode = $false
